$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(1726, 1703, 1714, 1730, 1727, 1724, 1708, 1718, 1710, 1722, 1712, 1706, 1733, 1735)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

$ws.Range("F15").Select()
